$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the gross expenditures and total maintenance inputs; the dependent
# formulas (D6, D7, D8, E8, D9, E9) recalc automatically.
$ws.Range("D3").Value = 413691.9
$ws.Range("D5").Value = 69691.58

# Move the active selection to the notes box on the right (G7:I16).
$ws.Range("G7:I16").Select()
